$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values for the new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered) from an existing
# header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-20
$iValues = @(11, 3, 9, 2, 4, 8, 2, 7, 3, 8, 3, 6, 7, 8, 5, 9, 9, 3, 5)
$jValues = @(11, 5, 9, 4, 5, 8, 5, 8, 5, 8, 5, 6, 8, 8, 5, 9, 9, 3, 5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
